# Generate Report for Handoff
# Updates the localization-status report: flips the "In Translation" rows to
# "Ready for handoff" and refreshes the associated handoff timestamps, then
# widens the (now longer) status columns to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Handoff generation timestamps ---
$wsOverview.Range("G2").Value = "2016-09-07 02:48:56"
$wsDeDe.Range("H2").Value     = "2016-09-07 02:48:56"
$wsZhCn.Range("H2").Value     = "2016-09-07 02:48:51"

# --- Widen the status columns so the longer text fits ---
# (new width matches the regenerated report's column sizing for the
# "Ready for handoff" / "Ready for handoff" status text)
$newStatusColWidth = 16.38265482584637
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth
$wsZhCn.Columns.Item(3).ColumnWidth     = $newStatusColWidth
$wsDeDe.Columns.Item(3).ColumnWidth     = $newStatusColWidth
